# Corretto data proc invasive, UOC neurochirurgia come reparto di riferimento
# Update computed statistics (Eff_cost, c_low, c_high, Eff_ln, ln_low, ln_high)
# for rows 2-14 on the active worksheet after recalculating with neurochirurgia
# (UOC neurochirurgia) set as the reference department.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 5715.21259259259
$ws.Range("D2").Value = 4755.59365903528
$ws.Range("E2").Value = 6674.8315261499
$ws.Range("F2").Value = 58.4736690232438
$ws.Range("G2").Value = 48.0412061479761
$ws.Range("H2").Value = 69.6413074923597
$ws.Range("C3").Value = 7591.30215827338
$ws.Range("D3").Value = 4180.26307059955
$ws.Range("E3").Value = 11002.3412459472
$ws.Range("F3").Value = 48.9998266489516
$ws.Range("G3").Value = 22.640584900941
$ws.Range("H3").Value = 81.0244818984657
$ws.Range("C4").Value = 6976.44827586207
$ws.Range("D4").Value = 4227.32196803053
$ws.Range("E4").Value = 9725.57458369361
$ws.Range("F4").Value = 58.9461068505946
$ws.Range("G4").Value = 34.0551918863844
$ws.Range("H4").Value = 88.4586827817343
$ws.Range("C5").Value = 6559.02857142857
$ws.Range("D5").Value = -492.471476130318
$ws.Range("E5").Value = 13610.5286189875
$ws.Range("F5").Value = 46.298727633075
$ws.Range("G5").Value = -4.67402246824763
$ws.Range("H5").Value = 124.527649873062
$ws.Range("C6").Value = 4728.76862745098
$ws.Range("D6").Value = 2969.80577555462
$ws.Range("E6").Value = 6487.73147934734
$ws.Range("F6").Value = 41.6683035115908
$ws.Range("G6").Value = 27.3690084481607
$ws.Range("H6").Value = 57.5729328851664
$ws.Range("C7").Value = 4012.29965156794
$ws.Range("D7").Value = 1549.7222994326
$ws.Range("E7").Value = 6474.87700370329
$ws.Range("F7").Value = 23.8659014613487
$ws.Range("G7").Value = 6.84649947783738
$ws.Range("H7").Value = 43.5962958057882
$ws.Range("C8").Value = 5905.4402173913
$ws.Range("D8").Value = 2889.42549041559
$ws.Range("E8").Value = 8921.45494436702
$ws.Range("F8").Value = 44.8845316038222
$ws.Range("G8").Value = 21.41121023058
$ws.Range("H8").Value = 72.8961226742783
$ws.Range("C9").Value = 5941.61733615222
$ws.Range("D9").Value = 4306.84686242448
$ws.Range("E9").Value = 7576.38780987996
$ws.Range("F9").Value = 54.4490455677735
$ws.Range("G9").Value = 38.4040986642213
$ws.Range("H9").Value = 72.3540553135569
$ws.Range("C10").Value = 7767.69244935543
$ws.Range("D10").Value = 6035.88844322261
$ws.Range("E10").Value = 9499.49645548826
$ws.Range("F10").Value = 63.5287467489932
$ws.Range("G10").Value = 46.6978910518555
$ws.Range("H10").Value = 82.2906302302846
$ws.Range("C11").Value = 9986.43418467583
$ws.Range("D11").Value = 7968.512908009
$ws.Range("E11").Value = 12004.3554613427
$ws.Range("F11").Value = 84.629022437725
$ws.Range("G11").Value = 63.4141853892643
$ws.Range("H11").Value = 108.598022534642
$ws.Range("C12").Value = 4695.12058823529
$ws.Range("D12").Value = 3482.82323385664
$ws.Range("E12").Value = 5907.41794261395
$ws.Range("F12").Value = 41.2873374027604
$ws.Range("G12").Value = 29.1788613355816
$ws.Range("H12").Value = 54.5307916788627
$ws.Range("C13").Value = 9635.97256097561
$ws.Range("D13").Value = 7277.18306734735
$ws.Range("E13").Value = 11994.7620546039
$ws.Range("F13").Value = 86.141824160658
$ws.Range("G13").Value = 60.8955424557884
$ws.Range("H13").Value = 115.349525369096
$ws.Range("C14").Value = 5363.75113122172
$ws.Range("D14").Value = 2704.66743607948
$ws.Range("E14").Value = 8022.83482636396
$ws.Range("F14").Value = 50.8517269867928
$ws.Range("G14").Value = 29.0932322739088
$ws.Range("H14").Value = 76.2775873998869
